# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Packham's Triumph, 17-kilo box) ahead of
# the existing Abate Fettel / 16-kilo records, pushing the old rows 263-287
# down to 265-289.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 263-264; everything from the old row 263 onward
# shifts down by two rows (old 263 -> 265, ... old 287 -> 289).
$ws.Rows("263:264").Insert()

# New row 263: Packham's Triumph, Primera, 17 kilos
$ws.Cells.Item(263,1).Value  = 11
$ws.Cells.Item(263,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(263,3).Value  = "Bíobío"
$ws.Cells.Item(263,4).Value  = 44491
$ws.Cells.Item(263,5).Value  = 8
$ws.Cells.Item(263,6).Value  = "Fruta"
$ws.Cells.Item(263,7).Value  = 100104
$ws.Cells.Item(263,8).Value  = "Frutos de pepita"
$ws.Cells.Item(263,9).Value  = 100104005
$ws.Cells.Item(263,10).Value = "Pera"
$ws.Cells.Item(263,11).Value = "Packham's Triumph"
$ws.Cells.Item(263,12).Value = "Primera"
$ws.Cells.Item(263,13).Value = 270
$ws.Cells.Item(263,14).Value = 8000
$ws.Cells.Item(263,15).Value = 9000
$ws.Cells.Item(263,16).Value = 8444
$ws.Cells.Item(263,17).Value = "$/caja 17 kilos empedrada"
$ws.Cells.Item(263,18).Value = "Región de O'Higgins"
$ws.Cells.Item(263,19).Value = 497
$ws.Cells.Item(263,20).Value = 17

# New row 264: Packham's Triumph, Segunda, 17 kilos
$ws.Cells.Item(264,1).Value  = 11
$ws.Cells.Item(264,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(264,3).Value  = "Bíobío"
$ws.Cells.Item(264,4).Value  = 44491
$ws.Cells.Item(264,5).Value  = 8
$ws.Cells.Item(264,6).Value  = "Fruta"
$ws.Cells.Item(264,7).Value  = 100104
$ws.Cells.Item(264,8).Value  = "Frutos de pepita"
$ws.Cells.Item(264,9).Value  = 100104005
$ws.Cells.Item(264,10).Value = "Pera"
$ws.Cells.Item(264,11).Value = "Packham's Triumph"
$ws.Cells.Item(264,12).Value = "Segunda"
$ws.Cells.Item(264,13).Value = 250
$ws.Cells.Item(264,14).Value = 7500
$ws.Cells.Item(264,15).Value = 7500
$ws.Cells.Item(264,16).Value = 7500
$ws.Cells.Item(264,17).Value = "$/caja 17 kilos empedrada"
$ws.Cells.Item(264,18).Value = "Región de O'Higgins"
$ws.Cells.Item(264,19).Value = 441
$ws.Cells.Item(264,20).Value = 17
